$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 2 4 '256.17'
Set-TextValue 2 5 '-0.20%'
Set-TextValue 3 4 '26.80'
Set-TextValue 3 5 '-0.14%'
Set-TextValue 4 4 '4.700'
Set-TextValue 4 5 '-0.68%'
Set-TextValue 5 4 '0.05941'
Set-TextValue 5 5 '0.37%'
Set-TextValue 6 4 '6.613'
Set-TextValue 6 5 '-0.75%'
Set-TextValue 7 4 '0.8512'
Set-TextValue 7 5 '-1.81%'
Set-TextValue 8 4 '0.9117'
Set-TextValue 8 5 '-3.27%'
Set-TextValue 9 4 '0.1380'
Set-TextValue 9 5 '-1.53%'
Set-TextValue 10 4 '0.04471'
Set-TextValue 10 5 '19.26%'
Set-TextValue 11 4 '0.06999'
Set-TextValue 11 5 '-1.67%'
Set-TextValue 12 4 '0.03060'
Set-TextValue 12 5 '-3.23%'
Set-TextValue 13 4 '0.09088'
Set-TextValue 13 5 '-1.83%'
Set-TextValue 14 4 '0.001521'
Set-TextValue 14 5 '-1.29%'
Set-TextValue 15 4 '0.0006033'
Set-TextValue 15 5 '-94.23%'
Set-TextValue 16 4 '0.006031'
Set-TextValue 16 5 '-0.97%'
Set-TextValue 17 4 '3.469'
Set-TextValue 17 5 '-0.80%'
Set-TextValue 18 4 '3.159'
Set-TextValue 18 5 '-1.37%'
Set-TextValue 20 4 '0.3029'
Set-TextValue 20 5 '-3.61%'
Set-TextValue 21 5 '1.35%'
Set-TextValue 22 4 '3.875'
Set-TextValue 22 5 '1.71%'
Set-TextValue 23 5 '1.19%'
Set-TextValue 24 4 '0.001216'
Set-TextValue 24 5 '-0.36%'
Set-TextValue 25 4 '0.004765'
Set-TextValue 25 5 '11.10%'
Set-TextValue 26 4 '0.0001201'
Set-TextValue 26 5 '-29.72%'
Set-TextValue 27 5 '2.15%'
Set-TextValue 40 4 '0.03773'
Set-TextValue 40 5 '-1.37%'
$ws.Cells.Item(41, 2).Value = 'KickToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue 41 4 '0.006195'
Set-TextValue 41 5 '-0.01%'
$ws.Cells.Item(42, 2).Value = 'BKEXToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue 42 4 '0.1095'
Set-TextValue 42 5 '-0.62%'
Set-TextValue 43 4 '0.002202'
Set-TextValue 43 5 '1.54%'
Set-TextValue 44 4 '0.01408'
Set-TextValue 44 5 '26.51%'
Set-TextValue 45 4 '0.00005314'
Set-TextValue 45 5 '-3.31%'
Set-TextValue 46 4 '0.00000000751'
Set-TextValue 46 5 '0.14%'
Set-TextValue 47 4 '0.04402'
Set-TextValue 47 5 '-50.24%'
Set-TextValue 48 4 '0.2258'
Set-TextValue 48 5 '9,157.30%'
Set-TextValue 49 5 '0.14%'
Set-TextValue 50 5 '0.14%'
